$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text, even when they
# look numeric (e.g. "13.30", "434.00"), so significant trailing zeros
# survive. A bare .Value assignment lets Excel auto-detect such strings
# as numbers, silently dropping trailing zeros / reformatting. Force a
# Text number format while writing, then clear formatting again so the
# cell is left with no explicit style (matching the workbook, where
# these data cells carry no per-cell style).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.956.25'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.151.06'
$ws.Range('D3').ClearFormats()
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.43'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.51'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.96%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.143.89'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.149'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.38'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.466'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('E13').Value = '  -2.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.98'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.669.61'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.048.36'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.149.67'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '489.58'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.72'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.711'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.68'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.04'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.30'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.69%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.76'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.21'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.98'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.76'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.47%  '
$ws.Range('E32').Value = '  -5.62%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.65'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.20%  '
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.05'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.71'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0740'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.96'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -7.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '434.00'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.97%  '
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.31'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.937.58'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.02%  '
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('E46').Value = '  -5.85%  '
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.88'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.94%  '
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.32'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.13%  '
